$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update status text "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Update timestamps
$overview.Range("G2").Value = "2016-09-06 19:20:43"
$zhcn.Range("H2").Value = "2016-09-06 19:20:38"
$dede.Range("H2").Value = "2016-09-06 19:20:43"

# Narrow the wide status/datetime columns (target stored width ~17.216;
# the COM ColumnWidth setter quantizes to 1/6-character steps, so feed it
# the character-width value whose quantized result lands closest to that
# target: 16.3333 -> stored width 17.1666...)
$overview.Range("E1").ColumnWidth = 16.3333333333333
$overview.Range("F1").ColumnWidth = 16.3333333333333
$zhcn.Range("C1").ColumnWidth = 16.3333333333333
$dede.Range("C1").ColumnWidth = 16.3333333333333
